# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range('D2').Value = '67.930.44'
$ws.Range('E2').Value = '  +1.16%  '

# Row 3 (Ethereum)
$ws.Range('D3').Value = '3.518.46'
$ws.Range('E3').Value = '  +0.18%  '

# Row 4 (TetherUSD)
$ws.Range('E4').Value = '  +0.00%  '

# Row 5 (BNB)
$ws.Range('D5').Value = '''601.55'
$ws.Range('E5').Value = '  +0.98%  '

# Row 6 (Solana)
$ws.Range('D6').Value = '''181.19'
$ws.Range('E6').Value = '  +4.37%  '

# Row 7 (USDC)
$ws.Range('E7').Value = '  +0.01%  '

# Row 8 (LidoStakedEther)
$ws.Range('D8').Value = '3.517.94'
$ws.Range('E8').Value = '  +0.18%  '

# Row 9 (XRP)
$ws.Range('D9').Value = '''0.596'
$ws.Range('E9').Value = '  +0.19%  '

# Row 10 (Dogecoin)
$ws.Range('D10').Value = '''0.141'
$ws.Range('E10').Value = '  +6.71%  '

# Row 11 (Toncoin)
$ws.Range('D11').Value = '''7.16'
$ws.Range('E11').Value = '  -1.69%  '

# Row 12 (Cardano)
$ws.Range('D12').Value = '''0.439'
$ws.Range('E12').Value = '  +0.76%  '

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range('D13').Value = '4.125.68'
$ws.Range('E13').Value = '  +0.13%  '

# Row 14 (Avalanche)
$ws.Range('D14').Value = '''32.70'
$ws.Range('E14').Value = '  +12.12%  '

# Row 15 (TRON)
$ws.Range('E15').Value = '  +1.01%  '

# Row 16 (WrappedBTC)
$ws.Range('D16').Value = '67.915.38'
$ws.Range('E16').Value = '  +1.20%  '

# Row 17 (ShibaInu)
$ws.Range('E17').Value = '  +0.38%  '

# Row 18 (WrappedEther)
$ws.Range('D18').Value = '3.521.10'
$ws.Range('E18').Value = '  -0.44%  '

# Row 19 (Polkadot)
$ws.Range('E19').Value = '  +0.29%  '

# Row 20 (Chainlink)
$ws.Range('D20').Value = '''14.50'
$ws.Range('E20').Value = '  +2.18%  '

# Row 21 (BitcoinCash)
$ws.Range('D21').Value = '''401.07'
$ws.Range('E21').Value = '  +1.15%  '

# Row 22 (Uniswap)
$ws.Range('E22').Value = '  -0.88%  '

# Row 23 (Litecoin)
$ws.Range('D23').Value = '''73.77'
$ws.Range('E23').Value = '  +0.85%  '

# Row 24 (Polygon)
$ws.Range('E24').Value = '  +0.89%  '

# Row 25 (Dai)
$ws.Range('E25').Value = '  -0.12%  '

# Row 26 (LEO)
$ws.Range('D26').Value = '''5.72'
$ws.Range('E26').Value = '  +1.10%  '

# Row 27 (PEPE)
$ws.Range('E27').Value = '  +1.16%  '

# Row 28 (InternetComputer(DFINITY))
$ws.Range('E28').Value = '  +2.35%  '

# Row 29 (Kaspa)
$ws.Range('E29').Value = '  -2.18%  '

# Row 30 (Binance-PegBSC-USD)
$ws.Range('D30').Value = '''0.998'
$ws.Range('E30').Value = '  -0.07%  '

# Row 31 (NEARProtocol)
$ws.Range('D31').Value = '''6.27'
$ws.Range('E31').Value = '  -1.13%  '

# Row 32 (Fetch.AI)
$ws.Range('E32').Value = '  -0.46%  '

# Row 33 (PancakeSwap)
$ws.Range('E33').Value = '  +1.50%  '

# Row 34 (EthereumClassic)
$ws.Range('D34').Value = '''23.93'
$ws.Range('E34').Value = '  +0.19%  '

# Row 35 (Aptos)
$ws.Range('D35').Value = '''7.50'
$ws.Range('E35').Value = '  +1.56%  '

# Row 36 (USDe)
$ws.Range('E36').Value = '  +0.08%  '

# Row 37 (ImmutableX)
$ws.Range('E37').Value = '  -2.63%  '

# Row 38 (Monero)
$ws.Range('D38').Value = '''163.05'
$ws.Range('E38').Value = '  -0.44%  '

# Row 39 (Mantle)
$ws.Range('E39').Value = '  -0.26%  '

# Row 40 (Stacks)
$ws.Range('D40').Value = '''1.92'
$ws.Range('E40').Value = '  +0.27%  '

# Row 41 (dogwifhat)
$ws.Range('D41').Value = '''2.79'
$ws.Range('E41').Value = '  +7.86%  '

# Row 42 (RenderToken)
$ws.Range('E42').Value = '  -1.10%  '

# Row 43 (Filecoin)
$ws.Range('D43').Value = '''4.70'
$ws.Range('E43').Value = '  -0.24%  '

# Row 44 (Maker)
$ws.Range('D44').Value = '2.890.09'
$ws.Range('E44').Value = '  +2.45%  '

# Row 45 (EnergySwap)
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '''0.0737'
$ws.Range('E45').Value = '  -1.83%  '

# Row 46 (Hedera)
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''26.37'
$ws.Range('E46').Value = '  -0.54%  '

# Row 47 (InjectiveProtocol)
$ws.Range('D47').Value = '''26.97'
$ws.Range('E47').Value = '  -1.88%  '

# Row 48 (OKB)
$ws.Range('E48').Value = '  -0.98%  '

# Row 49 (Bittensor)
$ws.Range('D49').Value = '''347.85'
$ws.Range('E49').Value = '  +1.75%  '

# Row 50 (VeChain)
$ws.Range('E50').Value = '  -0.18%  '

# Row 51 (ONDO)
$ws.Range('E51').Value = '  -1.49%  '
